# "Generate Report for Handback" -- mark the two localized hand-off rows as
# handed back (in sync with en-US), stamp the handback datetimes, and add
# "Latest Target File" / "Latest Handback File" links for each row on the
# zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Flip the shared "Ready for handoff" status to the handback message.
#    The same shared string is used on the Overview summary sheet and on
#    both locale sheets, so update every cell that currently shows it.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: stamp the handback datetime and add the target /
#    handback file hyperlinks (same targets as the existing source-file
#    and handoff-file links, since the file came back in sync).
# ---------------------------------------------------------------------
$zhcn.Range("H2").Value = "2016-03-21 14:16:14"
$zhcn.Range("H3").Value = "2016-03-21 14:16:14"

$zhcn.Range("F2").Value = "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/350445dcca170110908c3edd17e0451fcdd4c4d4/e2e/98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md",
    "",
    "",
    "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md"
) | Out-Null

$zhcn.Range("G2").Value = "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.zh-cn.xlf"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7ede1fb3c6b774ed86fc62da3390b202f21e829/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.zh-cn.xlf",
    "",
    "",
    "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.zh-cn.xlf"
) | Out-Null

$zhcn.Range("F3").Value = "db983710-bbac-4b42-b41f-82a2695092b1.md"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/350445dcca170110908c3edd17e0451fcdd4c4d4/e2e/db983710-bbac-4b42-b41f-82a2695092b1.md",
    "",
    "",
    "db983710-bbac-4b42-b41f-82a2695092b1.md"
) | Out-Null

$zhcn.Range("G3").Value = "db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.zh-cn.xlf"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e7ede1fb3c6b774ed86fc62da3390b202f21e829/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.zh-cn.xlf",
    "",
    "",
    "db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback timestamp and
#    locale-specific handoff-file targets.
# ---------------------------------------------------------------------
$dede.Range("H2").Value = "2016-03-21 14:16:20"
$dede.Range("H3").Value = "2016-03-21 14:16:20"

$dede.Range("F2").Value = "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md"
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/350445dcca170110908c3edd17e0451fcdd4c4d4/e2e/98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md",
    "",
    "",
    "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.md"
) | Out-Null

$dede.Range("G2").Value = "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.de-de.xlf"
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aab95d767fba08c121cbfb5bd10a848acc8b6fab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.de-de.xlf",
    "",
    "",
    "98eeb3e1-55c5-4d4c-bdb9-1615eed0db81.673a659cc13e3db0667e1c80fed06b94fb60ff31.de-de.xlf"
) | Out-Null

$dede.Range("F3").Value = "db983710-bbac-4b42-b41f-82a2695092b1.md"
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/350445dcca170110908c3edd17e0451fcdd4c4d4/e2e/db983710-bbac-4b42-b41f-82a2695092b1.md",
    "",
    "",
    "db983710-bbac-4b42-b41f-82a2695092b1.md"
) | Out-Null

$dede.Range("G3").Value = "db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.de-de.xlf"
$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aab95d767fba08c121cbfb5bd10a848acc8b6fab/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.de-de.xlf",
    "",
    "",
    "db983710-bbac-4b42-b41f-82a2695092b1.fd6c19042df4a767df8afef7fa958678bdf11567.de-de.xlf"
) | Out-Null

Write-Host "Handback report generated."
